$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; existing rows 11-30 shift down to 12-31.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new price-record data.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = 44498
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 100114002
$ws.Cells.Item(11, 7).Value = "Camote"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 20
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 20000
$ws.Cells.Item(11, 14).Value = "`$/malla 20 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 1000
$ws.Cells.Item(11, 17).Value = 20
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of
# column D (style carried over from the row above on insert, but set it
# explicitly too in case the engine did not propagate it).
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
